$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new test data rows for the "recordVerification" Data Provider
$ws.Range("A6").Value = "recordVerification"
$ws.Range("A7").Value = "recordVerification"

# Add hyperlinked email addresses (this also creates the "Hyperlink" cell style)
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:111@gmail.com", "", "", "111@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:beto@bob.com", "", "", "beto@bob.com")

$ws.Range("A7").Select()
